$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# "add cards for tires and cars": append a new time-tracking entry (row 15)
# mirroring the existing rows (Date in col A, Time [h] in col B).
# Copy the formatting (date number format/style) from the row above (A14)
# down onto the new date cell, then set the actual values.
$ws.Cells.Item(14, 1).Copy()
$ws.Cells.Item(15, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(15, 1).Value = 46048
$ws.Cells.Item(15, 2).Value = 1

# Excel moves the active selection to the next empty row after data entry
$ws.Range("B16").Select()
